$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2319.4666
$ws.Range("J19").Value = 1604.1111
$ws.Range("L19").Value = 1604.1111
$ws.Range("N19").Value = -1954.1111
$ws.Range("H70").Value = 1217.1
$ws.Range("I70").Value = 726.3333
$ws.Range("J70").Value = 1427.4286
$ws.Range("K70").Value = 2178.9999
$ws.Range("L70").Value = 4282.2858
$ws.Range("M70").Value = -1908.9999
$ws.Range("N70").Value = -4822.2858
$ws.Range("H73").Value = 1217.1
$ws.Range("I73").Value = 726.3333
$ws.Range("J73").Value = 1427.4286
$ws.Range("K73").Value = 2178.9999
$ws.Range("L73").Value = 4282.2858
$ws.Range("M73").Value = -1242.9999
$ws.Range("N73").Value = -6154.2858
$ws.Range("H96").Value = 222693.78
$ws.Range("I96").Value = 250405.5
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 751216.5
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -749843.5
$ws.Range("H103").Value = 381.96
$ws.Range("J103").Value = 529.1667
$ws.Range("L103").Value = 1587.5001
$ws.Range("N103").Value = -2759.5001
$ws.Range("H108").Value = 39683
$ws.Range("J108").Value = 39683
$ws.Range("L108").Value = 39683
$ws.Range("N108").Value = -47363
$ws.Range("H115").Value = 3932.1667
$ws.Range("I115").Value = 4518.6
$ws.Range("K115").Value = 13555.8
$ws.Range("M115").Value = -11988.8
$ws.Range("H132").Value = 2801.6453
$ws.Range("I132").Value = 1906.9231
$ws.Range("K132").Value = 5720.7693
$ws.Range("M132").Value = -3190.7693
$ws.Range("H134").Value = 99899.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99899.336
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = -110039.336
$ws.Range("H137").Value = 54626.367
$ws.Range("I137").Value = 1927
$ws.Range("J137").Value = 252249
$ws.Range("K137").Value = 5781
$ws.Range("L137").Value = 756747
$ws.Range("M137").Value = -3231
$ws.Range("N137").Value = -761847
$ws.Range("H138").Value = 1768.8695
$ws.Range("I138").Value = 1152.6666
$ws.Range("J138").Value = 2644.5264
$ws.Range("K138").Value = 3457.9998
$ws.Range("L138").Value = 7933.5792
$ws.Range("M138").Value = 1682.0002
$ws.Range("N138").Value = -18213.5792
$ws.Range("H141").Value = 183589
$ws.Range("I141").Value = 302981.66
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 908944.98
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -903764.98
$ws.Range("N141").Value = -23860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2630.353
$ws.Range("I63").Value = 2247.2727
$ws.Range("K63").Value = 2247.2727
$ws.Range("M63").Value = -1561.2727
$ws.Range("H66").Value = 2630.353
$ws.Range("I66").Value = 2247.2727
$ws.Range("K66").Value = 11236.3635
$ws.Range("M66").Value = -7804.363499999999
$ws.Range("H122").Value = 44002.6
$ws.Range("I122").Value = 44002.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 132007.8
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5697.56
$ws.Range("I132").Value = 6761.973
$ws.Range("J132").Value = 2668.077
$ws.Range("K132").Value = 20285.919
$ws.Range("L132").Value = 8004.231000000001
$ws.Range("M132").Value = -17755.919
$ws.Range("N132").Value = -13064.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 512.2857
$ws.Range("I22").Value = 516.2727
$ws.Range("J22").Value = 497.66666
$ws.Range("K22").Value = 516.2727
$ws.Range("L22").Value = 497.66666
$ws.Range("M22").Value = -343.2727
$ws.Range("N22").Value = -843.66666
$ws.Range("H58").Value = 117995
$ws.Range("J58").Value = 117995
$ws.Range("L58").Value = 117995
$ws.Range("H107").Value = 12432.9375
$ws.Range("J107").Value = 50333.332
$ws.Range("L107").Value = 50333.332
$ws.Range("N107").Value = -54173.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2371.122
$ws.Range("J31").Value = 2822.818
$ws.Range("L31").Value = 2822.818
$ws.Range("N31").Value = -3412.818
$ws.Range("H34").Value = 2371.122
$ws.Range("J34").Value = 2822.818
$ws.Range("L34").Value = 2822.818
$ws.Range("N34").Value = -3226.818
$ws.Range("H100").Value = 67890
$ws.Range("J100").Value = 67890
$ws.Range("L100").Value = 67890
$ws.Range("N100").Value = -70054
$ws.Range("H105").Value = 6943.64
$ws.Range("I105").Value = 5806.5713
$ws.Range("J105").Value = 8390.817999999999
$ws.Range("K105").Value = 5806.5713
$ws.Range("L105").Value = 8390.817999999999
$ws.Range("M105").Value = -4059.5713
$ws.Range("N105").Value = -11884.818
$ws.Range("H122").Value = 2695.7917
$ws.Range("I122").Value = 2538.5293
$ws.Range("K122").Value = 7615.5879
$ws.Range("M122").Value = -5165.5879
$ws.Range("H132").Value = 2027.725
$ws.Range("I132").Value = 2027.725
$ws.Range("K132").Value = 6083.174999999999
$ws.Range("M132").Value = -3553.174999999999
$ws.Range("H138").Value = 57780
$ws.Range("J138").Value = 57780
$ws.Range("L138").Value = 57780
$ws.Range("N138").Value = -68060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("H98").Value = 765.3333
$ws.Range("J98").Value = 765.3333
$ws.Range("L98").Value = 2295.9999
$ws.Range("N98").Value = -5291.9999
$ws.Range("H114").Value = 25001256
$ws.Range("I114").Value = 40000650
$ws.Range("J114").Value = 2266.3333
$ws.Range("K114").Value = 120001950
$ws.Range("L114").Value = 6798.999899999999
$ws.Range("M114").Value = -119998696
$ws.Range("N114").Value = -13306.9999
$ws.Range("H134").Value = 5575.125
$ws.Range("I134").Value = 1746.091
$ws.Range("K134").Value = 5238.272999999999
$ws.Range("M134").Value = -168.2729999999992
$ws.Range("H136").Value = 4659
$ws.Range("I136").Value = 1664.75
$ws.Range("K136").Value = 4994.25
$ws.Range("M136").Value = 105.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10735643
$ws.Range("I11").Value = 17107376
$ws.Range("K11").Value = 17107376
$ws.Range("M11").Value = -17107237
$ws.Range("H34").Value = 25172.4
$ws.Range("J34").Value = 25172.4
$ws.Range("L34").Value = 25172.4
$ws.Range("N34").Value = -25708.4
$ws.Range("H76").Value = 25172.4
$ws.Range("J76").Value = 25172.4
$ws.Range("L76").Value = 25172.4
$ws.Range("N76").Value = -25802.4
$ws.Range("H79").Value = 25172.4
$ws.Range("J79").Value = 25172.4
$ws.Range("L79").Value = 25172.4
$ws.Range("N79").Value = -27356.4
$ws.Range("H87").Value = 23353.334
$ws.Range("J87").Value = 23353.334
$ws.Range("L87").Value = 23353.334
$ws.Range("H90").Value = 23353.334
$ws.Range("J90").Value = 23353.334
$ws.Range("L90").Value = 70060.00199999999
$ws.Range("H126").Value = 3023.6667
$ws.Range("I126").Value = 2883.2
$ws.Range("J126").Value = 3199.25
$ws.Range("K126").Value = 8649.599999999999
$ws.Range("L126").Value = 9597.75
$ws.Range("M126").Value = -6179.599999999999
$ws.Range("N126").Value = -14537.75
$ws.Range("H132").Value = 2101.7222
$ws.Range("I132").Value = 2341.6428
$ws.Range("J132").Value = 1262
$ws.Range("K132").Value = 7024.928400000001
$ws.Range("L132").Value = 3786
$ws.Range("M132").Value = -4494.928400000001
$ws.Range("N132").Value = -8846

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8057.2354
$ws.Range("I7").Value = 7449
$ws.Range("J7").Value = 8926.143
$ws.Range("K7").Value = 7449
$ws.Range("L7").Value = 8926.143
$ws.Range("M7").Value = -7337
$ws.Range("N7").Value = -9150.143
$ws.Range("H16").Value = 1450.3334
$ws.Range("I16").Value = 1550.3636
$ws.Range("J16").Value = 350
$ws.Range("K16").Value = 1550.3636
$ws.Range("L16").Value = 350
$ws.Range("M16").Value = -1380.3636
$ws.Range("H82").Value = 747.2727
$ws.Range("I82").Value = 624.75
$ws.Range("K82").Value = 624.75
$ws.Range("M82").Value = -263.75
$ws.Range("H85").Value = 747.2727
$ws.Range("I85").Value = 624.75
$ws.Range("K85").Value = 624.75
$ws.Range("M85").Value = 623.25
$ws.Range("H122").Value = 15136.111
$ws.Range("I122").Value = 52000
$ws.Range("J122").Value = 4603.5713
$ws.Range("K122").Value = 156000
$ws.Range("L122").Value = 13810.7139
$ws.Range("M122").Value = -153550
$ws.Range("N122").Value = -18710.7139
$ws.Range("H126").Value = 8057.2354
$ws.Range("I126").Value = 7449
$ws.Range("J126").Value = 8926.143
$ws.Range("K126").Value = 22347
$ws.Range("L126").Value = 26778.429
$ws.Range("M126").Value = -19877
$ws.Range("N126").Value = -31718.429
$ws.Range("H132").Value = 9536.6
$ws.Range("I132").Value = 12255.9
$ws.Range("J132").Value = 4098
$ws.Range("K132").Value = 36767.7
$ws.Range("L132").Value = 12294
$ws.Range("M132").Value = -34237.7
$ws.Range("N132").Value = -17354

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8269.675999999999
$ws.Range("I132").Value = 8724.406000000001
$ws.Range("J132").Value = 5359.4
$ws.Range("K132").Value = 26173.218
$ws.Range("L132").Value = 16078.2
$ws.Range("M132").Value = -23643.218
$ws.Range("N132").Value = -21138.2
